$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D10 and F10 (H10 formula will recalculate automatically)
$ws.Range("D10").Value = 15
$ws.Range("F10").Value = 1

# Update B12, D12, F12 (H12 formula will recalculate automatically)
$ws.Range("B12").Value = 9
$ws.Range("D12").Value = 15
$ws.Range("F12").Value = 1

# Update the active cell selection on the sheet view
$ws.Range("L14").Select()
